# Apply the updated cryptocurrency price/volume figures (cols D and E) for rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Price = "30.311.84"; Volume = "  -1.05%  " },
    @{ Row = 3; Price = "1.868.43"; Volume = "  -0.76%  " },
    @{ Row = 4; Price = $null; Volume = "  +0.03%  " },
    @{ Row = 5; Price = "243.55"; Volume = "  -2.50%  " },
    @{ Row = 6; Price = $null; Volume = "  +0.06%  " },
    @{ Row = 7; Price = "0.4715"; Volume = "  -0.78%  " },
    @{ Row = 8; Price = "0.2867"; Volume = "  -2.33%  " },
    @{ Row = 9; Price = "0.06469"; Volume = "  -0.97%  " },
    @{ Row = 10; Price = "21.63"; Volume = "  -1.62%  " },
    @{ Row = 11; Price = "0.07784"; Volume = "  +0.60%  " },
    @{ Row = 12; Price = "95.95"; Volume = "  -1.01%  " },
    @{ Row = 13; Price = "1.867.41"; Volume = "  -0.82%  " },
    @{ Row = 14; Price = "0.7139"; Volume = "  -3.65%  " },
    @{ Row = 15; Price = "5.112"; Volume = "  -3.00%  " },
    @{ Row = 16; Price = "283.22"; Volume = "  +3.04%  " },
    @{ Row = 17; Price = "30.308.31"; Volume = "  -1.54%  " },
    @{ Row = 18; Price = "12.96"; Volume = "  -1.74%  " },
    @{ Row = 19; Price = "1.0000"; Volume = "  -0.01%  " },
    @{ Row = 20; Price = "0.000007454"; Volume = "  -1.23%  " },
    @{ Row = 21; Price = "2.111.50"; Volume = "  -0.88%  " },
    @{ Row = 22; Price = "1.001"; Volume = "  +0.13%  " },
    @{ Row = 23; Price = "5.243"; Volume = "  -1.68%  " },
    @{ Row = 24; Price = "6.243"; Volume = "  +0.12%  " },
    @{ Row = 25; Price = "162.60"; Volume = "  -0.79%  " },
    @{ Row = 26; Price = "8.969"; Volume = "  -2.74%  " },
    @{ Row = 27; Price = "18.68"; Volume = "  -1.03%  " },
    @{ Row = 28; Price = "1.877"; Volume = "  -2.19%  " },
    @{ Row = 29; Price = "0.09619"; Volume = "  -0.80%  " },
    @{ Row = 30; Price = $null; Volume = "  -2.10%  " },
    @{ Row = 31; Price = "1.481"; Volume = "  -1.96%  " },
    @{ Row = 32; Price = "4.195"; Volume = "  -2.40%  " },
    @{ Row = 33; Price = "4.111"; Volume = "  -1.02%  " },
    @{ Row = 34; Price = "0.04808"; Volume = "  -1.35%  " },
    @{ Row = 35; Price = "1.115"; Volume = "  -1.27%  " },
    @{ Row = 36; Price = "0.6842"; Volume = "  -2.24%  " },
    @{ Row = 37; Price = "2.710"; Volume = "  -0.31%  " },
    @{ Row = 38; Price = "0.01878"; Volume = "  -1.19%  " },
    @{ Row = 39; Price = "2.844"; Volume = "  +2.64%  " },
    @{ Row = 40; Price = "75.24"; Volume = "  +0.51%  " },
    @{ Row = 41; Price = "6.216"; Volume = "  -1.70%  " },
    @{ Row = 42; Price = "1.916"; Volume = "  -4.97%  " },
    @{ Row = 43; Price = "0.4177"; Volume = "  -1.60%  " },
    @{ Row = 44; Price = "0.9996"; Volume = "  -0.03%  " },
    @{ Row = 45; Price = "0.8231"; Volume = "  -2.07%  " },
    @{ Row = 46; Price = "100.61"; Volume = "  -2.13%  " },
    @{ Row = 47; Price = "9.685"; Volume = "  +3.56%  " },
    @{ Row = 48; Price = "6.980"; Volume = "  -1.09%  " },
    @{ Row = 49; Price = "35.01"; Volume = "  -1.76%  " },
    @{ Row = 50; Price = "892.59"; Volume = "  -2.80%  " },
    @{ Row = 51; Price = "0.05740"; Volume = "  +0.12%  " }
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($null -ne $u.Price) {
        $priceCell = $ws.Cells.Item($row, 4)   # column D
        # Force text storage so values such as "30.311.84" or "1.0000"
        # are not re-interpreted/rounded as numbers by Excel.
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $u.Price
        $priceCell.Style = "Normal"
    }
    $volCell = $ws.Cells.Item($row, 5)   # column E
    $volCell.Value = $u.Volume
}

Write-Host "Updated $($updates.Count) rows of crypto data."
